$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# C10 was previously synced in as text ("39"); normalize it to a real number.
$ws.Range("C10").Value = 39

# Append the new submission row synced on 2026-02-09 13:02:40.
$ws.Range("A11").Value = "2026-02-09 13:02:40"
$ws.Range("B11").Value = "Zainab Tijjani"

# Admission No keeps coming through as text from the sync, even though it
# looks numeric - force text storage so "38" stays a string, not a number.
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "38"
$ws.Range("C11").Style = "Normal"

$ws.Range("D11").Value = 7
